$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1440.3334
$ws.Range("I32").Value = 1299.75
$ws.Range("J32").Value = 1510.625
$ws.Range("K32").Value = 1299.75
$ws.Range("L32").Value = 1510.625
$ws.Range("M32").Value = -973.75
$ws.Range("N32").Value = -2162.625
$ws.Range("H33").Value = 306.45456
$ws.Range("I33").Value = 187.75
$ws.Range("J33").Value = 623
$ws.Range("K33").Value = 187.75
$ws.Range("L33").Value = 623
$ws.Range("M33").Value = 41.25
$ws.Range("N33").Value = -1081
$ws.Range("H64").Value = 7607.0713
$ws.Range("J64").Value = 9687.375
$ws.Range("L64").Value = 9687.375
$ws.Range("N64").Value = -10183.375
$ws.Range("H67").Value = 7607.0713
$ws.Range("J67").Value = 9687.375
$ws.Range("L67").Value = 9687.375
$ws.Range("N67").Value = -11403.375
$ws.Range("H86").Value = 3902.5334
$ws.Range("I86").Value = 2428.5715
$ws.Range("J86").Value = 5192.25
$ws.Range("K86").Value = 2428.5715
$ws.Range("L86").Value = 5192.25
$ws.Range("M86").Value = -1305.5715
$ws.Range("N86").Value = -7438.25
$ws.Range("H89").Value = 3902.5334
$ws.Range("I89").Value = 2428.5715
$ws.Range("J89").Value = 5192.25
$ws.Range("K89").Value = 12142.8575
$ws.Range("L89").Value = 25961.25
$ws.Range("M89").Value = -6526.8575
$ws.Range("N89").Value = -37193.25
$ws.Range("H137").Value = 836.55554
$ws.Range("I137").Value = 788.375
$ws.Range("K137").Value = 2365.125
$ws.Range("M137").Value = 184.875
$ws.Range("H138").Value = 3199.4
$ws.Range("J138").Value = 3487.7646
$ws.Range("L138").Value = 10463.2938
$ws.Range("N138").Value = -20743.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3078.4443
$ws.Range("I32").Value = 2820.923
$ws.Range("K32").Value = 2820.923
$ws.Range("M32").Value = -2533.923
$ws.Range("H45").Value = 3457.1667
$ws.Range("I45").Value = 1447
$ws.Range("J45").Value = 3859.2
$ws.Range("K45").Value = 1447
$ws.Range("L45").Value = 3859.2
$ws.Range("M45").Value = -1070
$ws.Range("N45").Value = -4613.2
$ws.Range("H88").Value = 1911.1538
$ws.Range("I88").Value = 2449.3333
$ws.Range("J88").Value = 1749.7
$ws.Range("K88").Value = 2449.3333
$ws.Range("L88").Value = 1749.7
$ws.Range("M88").Value = -2043.3333
$ws.Range("N88").Value = -2561.7
$ws.Range("H91").Value = 1911.1538
$ws.Range("I91").Value = 2449.3333
$ws.Range("J91").Value = 1749.7
$ws.Range("K91").Value = 2449.3333
$ws.Range("L91").Value = 1749.7
$ws.Range("M91").Value = -1045.3333
$ws.Range("N91").Value = -4557.7
$ws.Range("H102").Value = 2400
$ws.Range("J102").Value = 1500
$ws.Range("L102").Value = 1500
$ws.Range("N102").Value = -4744
$ws.Range("H122").Value = 1266.5
$ws.Range("I122").Value = 1120
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 3360
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -910
$ws.Range("N122").Value = -10897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2281.25
$ws.Range("I94").Value = 1510.2
$ws.Range("K94").Value = 1510.2
$ws.Range("M94").Value = -1059.2
$ws.Range("H99").Value = 2912.5
$ws.Range("I99").Value = 2912.5
$ws.Range("K99").Value = 2912.5
$ws.Range("M99").Value = -1414.5
$ws.Range("H105").Value = 3482.3333
$ws.Range("I105").Value = 3228.8
$ws.Range("J105").Value = 4750
$ws.Range("K105").Value = 3228.8
$ws.Range("L105").Value = 4750
$ws.Range("M105").Value = -1481.8
$ws.Range("N105").Value = -8244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2698.6
$ws.Range("I31").Value = 1998
$ws.Range("J31").Value = 3165.6667
$ws.Range("K31").Value = 1998
$ws.Range("L31").Value = 3165.6667
$ws.Range("M31").Value = -1703
$ws.Range("N31").Value = -3755.6667
$ws.Range("H34").Value = 2698.6
$ws.Range("I34").Value = 1998
$ws.Range("J34").Value = 3165.6667
$ws.Range("K34").Value = 1998
$ws.Range("L34").Value = 3165.6667
$ws.Range("M34").Value = -1796
$ws.Range("N34").Value = -3569.6667
$ws.Range("H94").Value = 2979.1428
$ws.Range("I94").Value = 3372.8572
$ws.Range("J94").Value = 2585.4285
$ws.Range("K94").Value = 3372.8572
$ws.Range("L94").Value = 2585.4285
$ws.Range("M94").Value = -2921.8572
$ws.Range("N94").Value = -3487.4285
$ws.Range("H105").Value = 15383.143
$ws.Range("J105").Value = 1990
$ws.Range("L105").Value = 1990
$ws.Range("N105").Value = -5484

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1730
$ws.Range("I69").Value = 598
$ws.Range("J69").Value = 3994
$ws.Range("K69").Value = 1794
$ws.Range("L69").Value = 11982
$ws.Range("M69").Value = -983
$ws.Range("N69").Value = -13604
$ws.Range("H72").Value = 1730
$ws.Range("I72").Value = 598
$ws.Range("J72").Value = 3994
$ws.Range("K72").Value = 5382
$ws.Range("L72").Value = 35946
$ws.Range("M72").Value = -1326
$ws.Range("N72").Value = -44058
$ws.Range("H98").Value = 276.14285
$ws.Range("I98").Value = 74.5
$ws.Range("J98").Value = 356.8
$ws.Range("K98").Value = 223.5
$ws.Range("L98").Value = 1070.4
$ws.Range("M98").Value = 1274.5
$ws.Range("N98").Value = -4066.4
$ws.Range("H107").Value = 1100.3334
$ws.Range("J107").Value = 1149
$ws.Range("L107").Value = 3447
$ws.Range("N107").Value = -7287
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -9930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3593.75
$ws.Range("I80").Value = 2687.5
$ws.Range("J80").Value = 4500
$ws.Range("K80").Value = 2687.5
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -1689.5
$ws.Range("N80").Value = -6496
$ws.Range("H83").Value = 3593.75
$ws.Range("I83").Value = 2687.5
$ws.Range("J83").Value = 4500
$ws.Range("K83").Value = 13437.5
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -8445.5
$ws.Range("N83").Value = -32484
$ws.Range("H102").Value = 2253
$ws.Range("I102").Value = 2139.6365
$ws.Range("K102").Value = 2139.6365
$ws.Range("M102").Value = -517.6365000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 559.26666
$ws.Range("I55").Value = 287
$ws.Range("J55").Value = 601.1539
$ws.Range("K55").Value = 287
$ws.Range("L55").Value = 601.1539
$ws.Range("M55").Value = -114
$ws.Range("N55").Value = -947.1539
$ws.Range("H68").Value = 2262.4
$ws.Range("I68").Value = 2262.4
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2262.4
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1513.4
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 2262.4
$ws.Range("I71").Value = 2262.4
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 11312
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -7568
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 3992.25
$ws.Range("I122").Value = 1989.6666
$ws.Range("K122").Value = 5968.9998
$ws.Range("M122").Value = -3518.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1070
$ws.Range("I122").Value = 1070
$ws.Range("K122").Value = 3210
$ws.Range("M122").Value = -760
